$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(112).Insert()

$ws.Cells.Item(112, 1).Value = 10
$ws.Cells.Item(112, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(112, 3).Value = 'La Araucanía'
$ws.Cells.Item(112, 4).Value = 44495
$ws.Cells.Item(112, 5).Value = 9
$ws.Cells.Item(112, 6).Value = 100114014
$ws.Cells.Item(112, 7).Value = 'Betarraga'
$ws.Cells.Item(112, 8).Value = 'Sin especificar'
$ws.Cells.Item(112, 9).Value = 'Primera'
$ws.Cells.Item(112, 10).Value = 50
$ws.Cells.Item(112, 11).Value = 10000
$ws.Cells.Item(112, 12).Value = 10000
$ws.Cells.Item(112, 13).Value = 10000
$ws.Cells.Item(112, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(112, 15).Value = 'Región del Maule'
$ws.Cells.Item(112, 16).Value = 833
$ws.Cells.Item(112, 17).Value = 12
$ws.Cells.Item(112, 18).Value = 'Hortaliza'
